# Update Handicap Index Course Handicap Report
# - H.I. / White / Blue columns (D/E/F) are stored as TEXT in this sheet
#   (e.g. "14.0" must keep its trailing zero), so values are written with
#   a leading apostrophe to force text, matching the existing data.
# - C23 holds the report's execution date/time as a normal date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$apos = [char]39

function Set-TextValue {
    param($addr, $val)
    $ws.Range($addr).Value = "$apos$val"
}

# Baker
Set-TextValue "D2" "11.9"

# Baumgarth
Set-TextValue "D4" "11.3"
Set-TextValue "E4" "11"

# Besse
Set-TextValue "D5" "14.0"
Set-TextValue "E5" "14"
Set-TextValue "F5" "17"

# Broyles
Set-TextValue "D6" "15.2"

# Carroll
Set-TextValue "D7" "9.3"

# Davis
Set-TextValue "D8" "14.4"
Set-TextValue "F8" "17"

# Fannon
Set-TextValue "D10" "13.0"

# Heard
Set-TextValue "D12" "5.9"
Set-TextValue "E12" "5"
Set-TextValue "F12" "7"

# Humphrey
Set-TextValue "D13" "18.2"

# Stewart
Set-TextValue "D15" "9.3"
Set-TextValue "E15" "9"
Set-TextValue "F15" "11"

# Traub
Set-TextValue "D17" "9.0"
Set-TextValue "F17" "11"

# Vela
Set-TextValue "D18" "18.0"
Set-TextValue "E18" "19"
Set-TextValue "F18" "21"

# Wickham
Set-TextValue "D19" "14.5"
Set-TextValue "F19" "17"

# Williams
Set-TextValue "D20" "16.3"
Set-TextValue "E20" "17"
Set-TextValue "F20" "19"

# Report execution date/time
$ws.Range("C23").Value = 43979.12501157408
